# Actualización completa del sistema de asistencia QR
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column order / labels ---
$ws.Cells.Item(1, 1).Value = "Fecha"
$ws.Cells.Item(1, 2).Value = "ID"
$ws.Cells.Item(1, 3).Value = "Nombre"
$ws.Cells.Item(1, 4).Value = "Departamento"
$ws.Cells.Item(1, 5).Value = "Entrada"
$ws.Cells.Item(1, 6).Value = "Salida"
$ws.Cells.Item(1, 7).Value = "Horas Trabajadas"
$ws.Cells.Item(1, 8).Value = "Dirección"

$direccion = "Institución Educativa No. 40657, Elías Aguirre, Selva Alegre, Chilina, Alto Selva Alegre, Arequipa, 04003, Perú"

# Column A holds dates written as plain dd/mm/yyyy text (not real date values).
# Temporarily mark the cells as Text so Excel doesn't auto-convert the string
# into a date serial number, then restore the default "Normal" style so no
# extra number formatting is left behind on the cell.
$ws.Range("A2:A4").NumberFormat = "@"

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value = "11/08/2025"
$ws.Cells.Item(2, 1).Style = "Normal"
$ws.Cells.Item(2, 2).Value = "E001"
$ws.Cells.Item(2, 3).Value = "Yadira Zeballos"
$ws.Cells.Item(2, 4).Value = "Tecnologia Digital "
$ws.Cells.Item(2, 5).Value = "09:34:51"
$ws.Cells.Item(2, 6).Value = "10:33:02"
$ws.Cells.Item(2, 7).Value = "58 minutos, 11 segundos"
$ws.Cells.Item(2, 8).Value = $direccion

# --- Row 3 (new) ---
$ws.Cells.Item(3, 1).Value = "11/08/2025"
$ws.Cells.Item(3, 1).Style = "Normal"
$ws.Cells.Item(3, 2).Value = "E002"
$ws.Cells.Item(3, 3).Value = "Melani Zeballos "
$ws.Cells.Item(3, 4).Value = "Marketing"
$ws.Cells.Item(3, 5).Value = "09:54:38"
$ws.Cells.Item(3, 6).Value = "No registrada"
$ws.Cells.Item(3, 7).Value = "N/A"
$ws.Cells.Item(3, 8).Value = $direccion

# --- Row 4 (new) ---
$ws.Cells.Item(4, 1).Value = "11/08/2025"
$ws.Cells.Item(4, 1).Style = "Normal"
$ws.Cells.Item(4, 2).Value = "E003"
$ws.Cells.Item(4, 3).Value = "Luis Lopez"
$ws.Cells.Item(4, 4).Value = "Tecnologia Digital "
$ws.Cells.Item(4, 5).Value = "10:33:13"
$ws.Cells.Item(4, 6).Value = "20:33:34"
$ws.Cells.Item(4, 7).Value = "10 horas, 21 segundos"
$ws.Cells.Item(4, 8).Value = $direccion
